# Update "想去人数" (column F) counts across all four sheets to reflect the
# latest scrape (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    2 = 142
    3 = 404
    5 = 42
    6 = 1259
    7 = 464
    9 = 217
    10 = 160
    11 = 191
    12 = 1065
    13 = 8
    14 = 277
    15 = 209
    16 = 1547
    17 = 570
    18 = 241
    19 = 362
    21 = 865
    22 = 1172
    25 = 2701
    26 = 1485
    28 = 59
    29 = 471
    31 = 1368
    33 = 1451
    34 = 170
    36 = 800
    37 = 667
    38 = 699
    39 = 891
    40 = 376
    41 = 266
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# ---- Sheet "演出" (performances) ----
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    15 = 674
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# ---- Sheet "本地生活" (local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$updates = @{
    2 = 868
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# ---- Sheet "全部类型" (all types, roll-up of the sheets above) ----
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    2 = 868
    3 = 142
    4 = 404
    6 = 42
    9 = 1259
    10 = 464
    12 = 217
    13 = 160
    14 = 191
    15 = 1065
    16 = 277
    18 = 209
    19 = 1547
    20 = 570
    21 = 241
    22 = 362
    25 = 1172
    26 = 2701
    28 = 1485
    31 = 59
    34 = 471
    36 = 1369
    40 = 1451
    41 = 800
    42 = 667
    43 = 699
    44 = 891
    45 = 376
    48 = 266
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
